# FiveGuys Deliverable 1 - Product/Sprint backlog update
# - Reword the "built in tutorial" user story (drop "some sort of")
# - Update Sprint # (col A) and Story Priority (col C) values for every
#   data row on the Sprint backlog sheet
# - Move the active selection to B29
# - Re-apply the sort over A3:G29 so the sortState range grows to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the wording of the tutorial user story (row 16, column B)
$ws.Range("B16").Value = "As a potential new user, I would like a built in tutorial, so that I can learn the game without having to look up a guide."

# 2) Updated Sprint # (A) and Story Priority (C) values for rows 3-29
$updates = @(
    @{Row=3;  A=1;  C=10},
    @{Row=4;  A=3;  C=3},
    @{Row=5;  A=3;  C=4},
    @{Row=6;  A=4;  C=5},
    @{Row=7;  A=5;  C=6},
    @{Row=8;  A=6;  C=5},
    @{Row=9;  A=7;  C=6},
    @{Row=10; A=8;  C=9},
    @{Row=11; A=9;  C=9},
    @{Row=12; A=10; C=6},
    @{Row=13; A=11; C=10},
    @{Row=14; A=12; C=4},
    @{Row=15; A=13; C=7},
    @{Row=16; A=14; C=5},
    @{Row=17; A=15; C=10},
    @{Row=18; A=16; C=9},
    @{Row=19; A=17; C=3},
    @{Row=20; A=18; C=7},
    @{Row=21; A=19; C=8},
    @{Row=22; A=20; C=7},
    @{Row=23; A=21; C=8},
    @{Row=24; A=22; C=9},
    @{Row=25; A=23; C=7},
    @{Row=26; A=24; C=9},
    @{Row=27; A=25; C=9},
    @{Row=28; A=26; C=9},
    @{Row=29; A=27; C=10}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}

# 3) Re-apply the sort across the full A3:G29 range (data is already in
#    ascending order by column A, so this is a no-op re-sort) so that the
#    saved sortState range grows from A4:D29 to A3:G29.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A3:A29"))
$ws.Sort.SetRange($ws.Range("A3:G29"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# 4) Move the selection to B29 (was C29)
$ws.Range("B29").Select()

# 5) Best-effort: restore the workbook window geometry recorded in the
#    saved file (not all hosts persist this, but set it for parity).
try {
    $win = $excel.ActiveWindow
    $win.Left = 18440
    $win.Top = -17250
    $win.Width = 19810
    $win.Height = 14610
} catch {}
